$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.76613666666667
$ws.Range("H2").Value = 53.29841
$ws.Range("I2").Value = 0.7872390387208499
$ws.Range("J2").Value = 0.7872390387208499
$ws.Range("M2").Value = 2.231113333333334
$ws.Range("N2").Value = 6.69334
$ws.Range("O2").Value = 0.01598125358798882
$ws.Range("P2").Value = 0.01598125358798882
$ws.Range("Q2").Value = 39.63826439882223
$ws.Range("R2").Value = 356.7443795894
$ws.Range("S2").Value = 0.01258106671216245
$ws.Range("T2").Value = 0.01258106671216245

$ws.Range("G3").Value = 17.76613666666667
$ws.Range("H3").Value = 53.29841
$ws.Range("I3").Value = 0.7872390387208499
$ws.Range("J3").Value = 0.7872390387208499
$ws.Range("O3").Value = 0.1634493267640196
$ws.Range("P3").Value = 0.1634493267640195
$ws.Range("Q3").Value = 405.4029675713967
$ws.Range("R3").Value = 3648.62670814257
$ws.Range("S3").Value = 0.1286736908812769
$ws.Range("T3").Value = 0.1286736908812768

$ws.Range("G4").Value = 17.76613666666667
$ws.Range("H4").Value = 53.29841
$ws.Range("I4").Value = 0.7872390387208499
$ws.Range("J4").Value = 0.7872390387208499
$ws.Range("M4").Value = 58.02175166666666
$ws.Range("N4").Value = 174.065255
$ws.Range("O4").Value = 0.4156043142904646
$ws.Range("P4").Value = 0.4156043142904646
$ws.Range("Q4").Value = 1030.822369749395
$ws.Range("R4").Value = 9277.40132774455
$ws.Range("S4").Value = 0.3271799408702633
$ws.Range("T4").Value = 0.3271799408702633

$ws.Range("G5").Value = 17.76613666666667
$ws.Range("H5").Value = 53.29841
$ws.Range("I5").Value = 0.7872390387208499
$ws.Range("J5").Value = 0.7872390387208499
$ws.Range("M5").Value = 15.16934033333333
$ws.Range("N5").Value = 45.508021
$ws.Range("O5").Value = 0.1086565487318021
$ws.Range("P5").Value = 0.1086565487318021
$ws.Range("Q5").Value = 269.5005735051789
$ws.Range("R5").Value = 2425.50516154661
$ws.Range("S5").Value = 0.08553867697434908
$ws.Range("T5").Value = 0.08553867697434908

$ws.Range("G6").Value = 17.76613666666667
$ws.Range("H6").Value = 53.29841
$ws.Range("I6").Value = 0.7872390387208499
$ws.Range("J6").Value = 0.7872390387208499
$ws.Range("M6").Value = 41.36709099999999
$ws.Range("N6").Value = 124.101273
$ws.Range("O6").Value = 0.2963085566257249
$ws.Range("P6").Value = 0.2963085566257249
$ws.Range("Q6").Value = 734.9333922084367
$ws.Range("R6").Value = 6614.40052987593
$ws.Range("S6").Value = 0.2332656632827982
$ws.Range("T6").Value = 0.2332656632827982

$ws.Range("I7").Value = 0.03648413815195897
$ws.Range("J7").Value = 0.03648413815195897
$ws.Range("M7").Value = 2.231113333333334
$ws.Range("N7").Value = 6.69334
$ws.Range("O7").Value = 0.01598125358798882
$ws.Range("P7").Value = 0.01598125358798882
$ws.Range("Q7").Value = 1.837012448951111
$ws.Range("R7").Value = 16.53311204056
$ws.Range("S7").Value = 0.000583062263745674
$ws.Range("T7").Value = 0.000583062263745674

$ws.Range("I8").Value = 0.03648413815195897
$ws.Range("J8").Value = 0.03648413815195897
$ws.Range("O8").Value = 0.1634493267640196
$ws.Range("P8").Value = 0.1634493267640195
$ws.Range("S8").Value = 0.005963307818503175
$ws.Range("T8").Value = 0.005963307818503174

$ws.Range("I9").Value = 0.03648413815195897
$ws.Range("J9").Value = 0.03648413815195897
$ws.Range("M9").Value = 58.02175166666666
$ws.Range("N9").Value = 174.065255
$ws.Range("O9").Value = 0.4156043142904646
$ws.Range("P9").Value = 0.4156043142904646
$ws.Range("Q9").Value = 47.77286681460221
$ws.Range("R9").Value = 429.9558013314199
$ws.Range("S9").Value = 0.01516296521912349
$ws.Range("T9").Value = 0.01516296521912349

$ws.Range("I10").Value = 0.03648413815195897
$ws.Range("J10").Value = 0.03648413815195897
$ws.Range("M10").Value = 15.16934033333333
$ws.Range("N10").Value = 45.508021
$ws.Range("O10").Value = 0.1086565487318021
$ws.Range("P10").Value = 0.1086565487318021
$ws.Range("Q10").Value = 12.48984828264044
$ws.Range("R10").Value = 112.408634543764
$ws.Range("S10").Value = 0.003964240535046131
$ws.Range("T10").Value = 0.003964240535046131

$ws.Range("I11").Value = 0.03648413815195897
$ws.Range("J11").Value = 0.03648413815195897
$ws.Range("M11").Value = 41.36709099999999
$ws.Range("N11").Value = 124.101273
$ws.Range("O11").Value = 0.2963085566257249
$ws.Range("P11").Value = 0.2963085566257249
$ws.Range("Q11").Value = 34.06006320188133
$ws.Range("R11").Value = 306.540568816932
$ws.Range("S11").Value = 0.01081056231554051
$ws.Range("T11").Value = 0.01081056231554051

$ws.Range("G12").Value = 3.885299333333334
$ws.Range("H12").Value = 11.655898
$ws.Range("I12").Value = 0.1721623203571791
$ws.Range("J12").Value = 0.172162320357179
$ws.Range("M12").Value = 2.231113333333334
$ws.Range("N12").Value = 6.69334
$ws.Range("O12").Value = 0.01598125358798882
$ws.Range("P12").Value = 0.01598125358798882
$ws.Range("Q12").Value = 8.668543146591112
$ws.Range("R12").Value = 78.01688831932
$ws.Range("S12").Value = 0.002751369699924648
$ws.Range("T12").Value = 0.002751369699924648

$ws.Range("G13").Value = 3.885299333333334
$ws.Range("H13").Value = 11.655898
$ws.Range("I13").Value = 0.1721623203571791
$ws.Range("J13").Value = 0.172162320357179
$ws.Range("O13").Value = 0.1634493267640196
$ws.Range("P13").Value = 0.1634493267640195
$ws.Range("Q13").Value = 88.65809766012734
$ws.Range("R13").Value = 797.922878941146
$ws.Range("S13").Value = 0.02813981535651238
$ws.Range("T13").Value = 0.02813981535651237

$ws.Range("G14").Value = 3.885299333333334
$ws.Range("H14").Value = 11.655898
$ws.Range("I14").Value = 0.1721623203571791
$ws.Range("J14").Value = 0.172162320357179
$ws.Range("M14").Value = 58.02175166666666
$ws.Range("N14").Value = 174.065255
$ws.Range("O14").Value = 0.4156043142904646
$ws.Range("P14").Value = 0.4156043142904646
$ws.Range("Q14").Value = 225.4318730693322
$ws.Range("R14").Value = 2028.88685762399
$ws.Range("S14").Value = 0.0715514030987007
$ws.Range("T14").Value = 0.07155140309870069

$ws.Range("G15").Value = 3.885299333333334
$ws.Range("H15").Value = 11.655898
$ws.Range("I15").Value = 0.1721623203571791
$ws.Range("J15").Value = 0.172162320357179
$ws.Range("M15").Value = 15.16934033333333
$ws.Range("N15").Value = 45.508021
$ws.Range("O15").Value = 0.1086565487318021
$ws.Range("P15").Value = 0.1086565487318021
$ws.Range("Q15").Value = 58.93742788420645
$ws.Range("R15").Value = 530.436850957858
$ws.Range("S15").Value = 0.01870656355166996
$ws.Range("T15").Value = 0.01870656355166996

$ws.Range("G16").Value = 3.885299333333334
$ws.Range("H16").Value = 11.655898
$ws.Range("I16").Value = 0.1721623203571791
$ws.Range("J16").Value = 0.172162320357179
$ws.Range("M16").Value = 41.36709099999999
$ws.Range("N16").Value = 124.101273
$ws.Range("O16").Value = 0.2963085566257249
$ws.Range("P16").Value = 0.2963085566257249
$ws.Range("Q16").Value = 160.7235310842393
$ws.Range("R16").Value = 1446.511779758154
$ws.Range("S16").Value = 0.05101316865037139
$ws.Range("T16").Value = 0.05101316865037138

$ws.Range("G17").Value = 0.09285466666666665
$ws.Range("H17").Value = 0.278564
$ws.Range("I17").Value = 0.004114502770011991
$ws.Range("J17").Value = 0.004114502770011991
$ws.Range("M17").Value = 2.231113333333334
$ws.Range("N17").Value = 6.69334
$ws.Range("O17").Value = 0.01598125358798882
$ws.Range("P17").Value = 0.01598125358798882
$ws.Range("Q17").Value = 0.2071692848622222
$ws.Range("R17").Value = 1.86452356376
$ws.Range("S17").Value = 0.00006575491215604406
$ws.Range("T17").Value = 0.00006575491215604406

$ws.Range("G18").Value = 0.09285466666666665
$ws.Range("H18").Value = 0.278564
$ws.Range("I18").Value = 0.004114502770011991
$ws.Range("J18").Value = 0.004114502770011991
$ws.Range("O18").Value = 0.1634493267640196
$ws.Range("P18").Value = 0.1634493267640195
$ws.Range("Q18").Value = 2.118837546158666
$ws.Range("R18").Value = 19.069537915428
$ws.Range("S18").Value = 0.0006725127077271536
$ws.Range("T18").Value = 0.0006725127077271535

$ws.Range("G19").Value = 0.09285466666666665
$ws.Range("H19").Value = 0.278564
$ws.Range("I19").Value = 0.004114502770011991
$ws.Range("J19").Value = 0.004114502770011991
$ws.Range("M19").Value = 58.02175166666666
$ws.Range("N19").Value = 174.065255
$ws.Range("O19").Value = 0.4156043142904646
$ws.Range("P19").Value = 0.4156043142904646
$ws.Range("Q19").Value = 5.387590410424443
$ws.Range("R19").Value = 48.48831369381999
$ws.Range("S19").Value = 0.001710005102377051
$ws.Range("T19").Value = 0.001710005102377051

$ws.Range("G20").Value = 0.09285466666666665
$ws.Range("H20").Value = 0.278564
$ws.Range("I20").Value = 0.004114502770011991
$ws.Range("J20").Value = 0.004114502770011991
$ws.Range("M20").Value = 15.16934033333333
$ws.Range("N20").Value = 45.508021
$ws.Range("O20").Value = 0.1086565487318021
$ws.Range("P20").Value = 0.1086565487318021
$ws.Range("Q20").Value = 1.408544040204889
$ws.Range("R20").Value = 12.676896361844
$ws.Range("S20").Value = 0.0004470676707369427
$ws.Range("T20").Value = 0.0004470676707369427

$ws.Range("G21").Value = 0.09285466666666665
$ws.Range("H21").Value = 0.278564
$ws.Range("I21").Value = 0.004114502770011991
$ws.Range("J21").Value = 0.004114502770011991
$ws.Range("M21").Value = 41.36709099999999
$ws.Range("N21").Value = 124.101273
$ws.Range("O21").Value = 0.2963085566257249
$ws.Range("P21").Value = 0.2963085566257249
$ws.Range("Q21").Value = 3.841127445774666
$ws.Range("R21").Value = 34.570147011972
$ws.Range("S21").Value = 0.0012191623770148
$ws.Range("T21").Value = 0.0012191623770148
